# Junction_Flooding_226.xlsx edit:
#  - apply "custom accuracy" (round to 2 decimal places) to the numeric
#    measurement cells of row 5 that actually need it
#  - drop the now-superfluous duplicate data row (row 6) -> "데이터 1000개"
#    trims the dataset; dimension collapses from A1:AH6 to A1:AH5
#  - column Q (17th column) narrows from width 8 to width 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round the row-5 cells whose value has more than 2 decimal digits ---
$cols = @(3, 4, 6, 8, 9, 10, 11, 13, 14, 16, 17, 19, 20, 21, 23, 24, 25, 26, 27, 29, 30, 31, 32, 33, 34)
foreach ($c in $cols) {
    $cell = $ws.Cells.Item(5, $c)
    $cell.Value = [Math]::Round($cell.Value(), 2)
}

# --- Remove row 6 entirely (shrinks used range to A1:AH5) ---
$ws.Rows.Item(6).Delete()

# --- Narrow column Q (column 17) from width 8 to width 7 ---
$ws.Columns.Item(17).ColumnWidth = 6.17
